$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '68.385.71'
$ws.Range('E2').Value = '  -1.09%  '
$ws.Range('D3').Value = '3.833.98'
$ws.Range('E3').Value = '  +2.46%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = "'599.18"
$ws.Range('E5').Value = '  -0.40%  '
$ws.Range('D6').Value = "'161.56"
$ws.Range('D7').Value = '3.833.08'
$ws.Range('E7').Value = '  +2.48%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = "'0.527"
$ws.Range('E9').Value = '  -2.52%  '
$ws.Range('E10').Value = '  -2.27%  '
$ws.Range('D11').Value = "'6.29"
$ws.Range('E11').Value = '  -1.35%  '
$ws.Range('D12').Value = "'0.456"
$ws.Range('E12').Value = '  -0.84%  '
$ws.Range('D13').Value = "'36.57"
$ws.Range('E13').Value = '  -3.80%  '
$ws.Range('E14').Value = '  -2.41%  '
$ws.Range('D15').Value = '4.485.10'
$ws.Range('E15').Value = '  +2.79%  '
$ws.Range('D16').Value = '3.840.85'
$ws.Range('D17').Value = '68.692.21'
$ws.Range('E17').Value = '  -0.52%  '
$ws.Range('D18').Value = "'7.50"
$ws.Range('E18').Value = '  +1.91%  '
$ws.Range('E19').Value = '  -0.40%  '
$ws.Range('E20').Value = '  -1.55%  '
$ws.Range('D21').Value = "'11.24"
$ws.Range('E21').Value = '  +0.81%  '
$ws.Range('D22').Value = "'482.88"
$ws.Range('E22').Value = '  -2.06%  '
$ws.Range('D23').Value = "'0.714"
$ws.Range('E23').Value = '  -1.85%  '
$ws.Range('E24').Value = '  +6.03%  '
$ws.Range('D25').Value = "'83.76"
$ws.Range('E25').Value = '  -1.46%  '
$ws.Range('E26').Value = '  -3.78%  '
$ws.Range('D27').Value = "'12.03"
$ws.Range('D28').Value = "'0.997"
$ws.Range('E28').Value = '  -0.24%  '
$ws.Range('D29').Value = "'9.89"
$ws.Range('E29').Value = '  -1.82%  '
$ws.Range('D30').Value = "'2.93"
$ws.Range('E30').Value = '  -1.39%  '
$ws.Range('D31').Value = '3.990.10'
$ws.Range('E31').Value = '  +2.62%  '
$ws.Range('D32').Value = "'7.80"
$ws.Range('E32').Value = '  -4.61%  '
$ws.Range('E33').Value = '  -4.57%  '
$ws.Range('D34').Value = "'31.86"
$ws.Range('E34').Value = '  +1.12%  '
$ws.Range('D35').Value = '3.785.18'
$ws.Range('E35').Value = '  +3.00%  '
$ws.Range('E36').Value = '  -2.17%  '
$ws.Range('E37').Value = '  +1.30%  '
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('D39').Value = "'5.83"
$ws.Range('E39').Value = '  -1.95%  '
$ws.Range('E40').Value = '  +0.20%  '
$ws.Range('E41').Value = '  -3.01%  '
$ws.Range('E42').Value = '  -3.67%  '
$ws.Range('D43').Value = "'425.78"
$ws.Range('E43').Value = '  +0.33%  '
$ws.Range('D44').Value = "'48.46"
$ws.Range('E44').Value = '  -0.71%  '
$ws.Range('E45').Value = '  -1.42%  '
$ws.Range('E46').Value = '  -0.01%  '
$ws.Range('D47').Value = "'8.33"
$ws.Range('E47').Value = '  -1.52%  '
$ws.Range('D48').Value = "'142.67"
$ws.Range('E48').Value = '  +0.37%  '
$ws.Range('D49').Value = '2.826.62'
$ws.Range('E49').Value = '  +1.62%  '
$ws.Range('D50').Value = "'25.85"
$ws.Range('E50').Value = '  +13.44%  '
$ws.Range('D51').Value = "'0.0354"
$ws.Range('E51').Value = '  +0.36%  '
